$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8484175205230713
$ws.Range("B1").Value = 2.923980474472046
$ws.Range("C1").Value = 4.633658885955811
$ws.Range("D1").Value = 2.794754981994629
$ws.Range("E1").Value = 1.435720443725586
